# Daily attendance processing - 2025-12-12 16:32:04
#
# For rows in column G ("Recorded By") that list a real user together with a
# "System"/backup entry, swap the first two entries in that comma-separated
# list (leaving any further entries, such as a trailing "System", in place).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

# Only the specific "Recorded By" combinations seen below should have their
# first two entries reordered; lists that are just a backup/admin account
# plus "System" (with no real user attached) are left as-is.
$targetValues = @(
    "dnasr281@gmail.com, System",
    "backup@backdoor.com, system, System",
    "dnasr281@gmail.com, admin@admin.com"
)

for ($row = 2; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)   # Column G = "Recorded By"
    $value = $cell.Value2

    if ($value -ne $null -and ($targetValues -contains $value)) {
        $parts = @($value -split ",\s*")
        if ($parts.Count -gt 2) {
            $tail = $parts[2..($parts.Count - 1)]
        } else {
            $tail = @()
        }
        $newParts = @($parts[1], $parts[0]) + $tail
        $newValue = [string]::Join(", ", $newParts)
        $cell.Value2 = $newValue
    }
}
